# Roll the IGCC Netting Flows history window forward by one day:
#  - every timestamp in column A (rows 2-201) advances by exactly 1 day
#  - the "Lookup" text in column E is rebuilt from the (unchanged) Quarter
#    number in column D, using the new date-of-month for each block
#  - columns B, C and D are left untouched
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 201

for ($r = 2; $r -le $lastRow; $r++) {
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value2 = $dateCell.Value2 + 1

    if ($r -le 105) {
        $day = 17
    } else {
        $day = 18
    }

    $quarter = $ws.Cells.Item($r, 4).Value2
    $lookupText = "{0:D2}.02.2026{1}" -f $day, $quarter
    $ws.Cells.Item($r, 5).Value = $lookupText
}
